# "adding averages and more checks"
#
# Training Dashboard sheet: the "PERIOD TO EXPIRE" (H) values all drop by 8
# days because "LAST UPDATE" (I) moved forward from 08-Sep-2025 to
# 16-Sep-2025 for every training row (3-22).
#
# Exam Dashboard sheet: the COMMENTS column (E) text changes from the
# generic "OK" to the more descriptive "date is valid" for rows 3-8, and
# that column is widened a bit (10 -> 15 chars) to fit the longer text.
#
# Both sheets also pick up a styling tweak: the bold header font (and, as a
# side effect of how the font was re-used, the dashboard title font too)
# becomes white so it reads cleanly against the dark-blue header fill.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Training Dashboard
# ---------------------------------------------------------------------
$wsTraining = $wb.Worksheets.Item("Training Dashboard")

$periodToExpire = @{
    3  = 484
    4  = 518
    5  = 484
    6  = 485
    7  = 483
    8  = 525
    9  = 484
    10 = 485
    11 = 485
    12 = 485
    13 = 484
    14 = 706
    15 = 485
    16 = 113
    17 = 169
    18 = 169
    19 = 155
    20 = 278
    21 = 313
    22 = 313
}

foreach ($row in 3..22) {
    $wsTraining.Range("H$row").Value = $periodToExpire[$row]
    # Force literal text (matches the original inline-string cell) instead
    # of letting the date-looking text auto-convert into a real date.
    $wsTraining.Range("I$row").Value = "'16-Sep-2025"
}

# Title (A1) and column header row (A2:K2) -> bold white text.
$wsTraining.Range("A1").Font.Size = 11
$wsTraining.Range("A1").Font.Color = 16777215
$wsTraining.Range("A2:K2").Font.Color = 16777215

# ---------------------------------------------------------------------
# Exam Dashboard
# ---------------------------------------------------------------------
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

foreach ($row in 3..8) {
    $wsExam.Range("E$row").Value = "date is valid"
}

# Widen the COMMENTS column (E) from 10 to 15 characters.
$wsExam.Columns.Item(5).ColumnWidth = 14.14

# Title (A1) and column header row (A2:G2) -> bold white text.
$wsExam.Range("A1").Font.Size = 11
$wsExam.Range("A1").Font.Color = 16777215
$wsExam.Range("A2:G2").Font.Color = 16777215
